$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store the given text verbatim, without Excel
# auto-converting numeric-looking strings into numbers, and without leaving
# a stray number-format style on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "68.390.35"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "3.751.65"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.22%  "
Set-TextValue $ws.Range("D5") "595.52"
$ws.Range("E5").Value = "  -0.27%  "
Set-TextValue $ws.Range("D6") "166.69"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("D7").Value = "3.746.75"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -1.29%  "
Set-TextValue $ws.Range("D13") "0.0000258"
$ws.Range("E13").Value = "  -6.68%  "
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "4.380.55"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "3.742.42"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "68.383.72"
$ws.Range("E17").Value = "  +1.00%  "
Set-TextValue $ws.Range("D18") "17.85"
$ws.Range("E18").Value = "  -4.53%  "
Set-TextValue $ws.Range("D19") "7.01"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("E21").Value = "  +1.15%  "
Set-TextValue $ws.Range("D22") "467.80"
$ws.Range("E22").Value = "  -0.17%  "
Set-TextValue $ws.Range("D23") "0.698"
$ws.Range("E23").Value = "  -3.09%  "
Set-TextValue $ws.Range("D24") "84.30"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("E26").Value = "  -0.89%  "
Set-TextValue $ws.Range("D27") "12.03"
$ws.Range("E27").Value = "  -1.07%  "
Set-TextValue $ws.Range("D28") "10.12"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "3.897.67"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  -4.69%  "
Set-TextValue $ws.Range("D32") "7.27"
$ws.Range("E32").Value = "  -5.03%  "
Set-TextValue $ws.Range("D33") "29.86"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -2.22%  "
Set-TextValue $ws.Range("D35") "9.24"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D37").Value = "3.706.30"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("E39").Value = "  -11.88%  "
$ws.Range("E40").Value = "  +0.48%  "
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.26%  "
Set-TextValue $ws.Range("D45") "0.304"
$ws.Range("E45").Value = "  -2.46%  "
Set-TextValue $ws.Range("D46") "8.59"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D47") "43.10"
$ws.Range("E47").Value = "  +10.28%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D48") "1.93"
$ws.Range("E48").Value = "  -0.80%  "
Set-TextValue $ws.Range("D49") "45.69"
$ws.Range("E49").Value = "  -0.29%  "
Set-TextValue $ws.Range("D50") "146.50"
$ws.Range("E50").Value = "  +4.89%  "
Set-TextValue $ws.Range("D51") "391.05"
$ws.Range("E51").Value = "  -1.61%  "
